$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.885218501091003
$ws.Range("B1").Value = 4.779308319091797
$ws.Range("C1").Value = 3.601540088653564
$ws.Range("D1").Value = 0.9005609750747681
$ws.Range("E1").Value = 0.4729504883289337
